# results_others.xlsx -- "print lines weggehaald bij de code"
#
# The "medium" sheet gained a new top data point (J1:K1), pushing the
# existing J:K column values down by one row (the former row 34 values
# land on the newly-used row 35). The other sheets only have their
# saved cursor/selection state changed (no data edits there).

$wb = $excel.ActiveWorkbook

# --- results_others: selection narrows from J1:K28 to J1:J28 --------------
$wsOthers = $wb.Worksheets.Item("results_others")
$wsOthers.Range("J1:J28").Select()

# --- easy: selection moves off the J:K block onto a single cell -----------
$wsEasy = $wb.Worksheets.Item("easy")
$wsEasy.Range("H23").Select()

# --- Hard: selection moves to A22:K22 (and it stops being the active tab) -
$wsHard = $wb.Worksheets.Item("Hard")
$wsHard.Range("A22:K22").Select()

# --- medium: shift the J:K results down one row and add a new top row ----
$wsMedium = $wb.Worksheets.Item("medium")

for ($row = 34; $row -ge 1; $row--) {
    $jValue = $wsMedium.Cells.Item($row, 10).Value2
    $kValue = $wsMedium.Cells.Item($row, 11).Value2
    $wsMedium.Cells.Item($row + 1, 10).Value = $jValue
    $wsMedium.Cells.Item($row + 1, 11).Value = $kValue
}

$wsMedium.Cells.Item(1, 10).Value = 140.35
$wsMedium.Cells.Item(1, 11).Value = 24.53

# medium becomes the active sheet/tab, selection anchored on the J column
$wsMedium.Activate()
$wsMedium.Range("J1:J35").Select()
